# Refresh cryptos list: updated price/volume figures for all 50 coin rows,
# plus two coin-pair reorders (Stellar/Monero swap at rows 26-27,
# Elrond/Decentraland swap at rows 49-50).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new values are plain numeric-looking strings (e.g. "0.9977") that Excel
# would otherwise auto-convert to numbers on assignment; the source data is
# text (t="inlineStr" in the original), so force Text format while writing,
# then clear the format again so the cells end up styled exactly as before.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "31.242.91"
$ws.Cells.Item(2, 5).Value = "  +2.92%  "
$ws.Cells.Item(3, 4).Value = "1.987.98"
$ws.Cells.Item(3, 5).Value = "  +6.28%  "
$ws.Cells.Item(4, 4).Value = "0.9977"
$ws.Cells.Item(4, 5).Value = "  -0.25%  "
$ws.Cells.Item(5, 4).Value = "0.8107"
$ws.Cells.Item(5, 5).Value = "  +72.16%  "
$ws.Cells.Item(6, 4).Value = "253.91"
$ws.Cells.Item(6, 5).Value = "  +4.13%  "
$ws.Cells.Item(7, 4).Value = "0.9979"
$ws.Cells.Item(7, 5).Value = "  -0.23%  "
$ws.Cells.Item(8, 4).Value = "0.3438"
$ws.Cells.Item(8, 5).Value = "  +19.45%  "
$ws.Cells.Item(9, 4).Value = "25.65"
$ws.Cells.Item(9, 5).Value = "  +16.62%  "
$ws.Cells.Item(10, 4).Value = "0.06972"
$ws.Cells.Item(10, 5).Value = "  +8.03%  "
$ws.Cells.Item(11, 4).Value = "0.8443"
$ws.Cells.Item(11, 5).Value = "  +16.51%  "
$ws.Cells.Item(12, 4).Value = "0.08110"
$ws.Cells.Item(13, 4).Value = "1.986.88"
$ws.Cells.Item(13, 5).Value = "  +6.22%  "
$ws.Cells.Item(14, 4).Value = "100.67"
$ws.Cells.Item(14, 5).Value = "  +4.66%  "
$ws.Cells.Item(15, 4).Value = "5.507"
$ws.Cells.Item(15, 5).Value = "  +7.39%  "
$ws.Cells.Item(16, 4).Value = "272.61"
$ws.Cells.Item(16, 5).Value = "  -2.26%  "
$ws.Cells.Item(17, 4).Value = "31.229.80"
$ws.Cells.Item(17, 5).Value = "  +2.90%  "
$ws.Cells.Item(18, 4).Value = "13.96"
$ws.Cells.Item(18, 5).Value = "  +7.39%  "
$ws.Cells.Item(19, 4).Value = "0.000007955"
$ws.Cells.Item(19, 5).Value = "  +6.07%  "
$ws.Cells.Item(20, 4).Value = "5.808"
$ws.Cells.Item(20, 5).Value = "  +10.83%  "
$ws.Cells.Item(21, 4).Value = "2.245.51"
$ws.Cells.Item(21, 5).Value = "  +6.40%  "
$ws.Cells.Item(22, 4).Value = "0.9990"
$ws.Cells.Item(22, 5).Value = "  -0.10%  "
$ws.Cells.Item(23, 4).Value = "0.9975"
$ws.Cells.Item(23, 5).Value = "  -0.28%  "
$ws.Cells.Item(24, 4).Value = "6.935"
$ws.Cells.Item(24, 5).Value = "  +11.30%  "
$ws.Cells.Item(25, 4).Value = "9.743"
$ws.Cells.Item(25, 5).Value = "  +7.69%  "
$ws.Cells.Item(26, 2).Value = "Monero"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(26, 4).Value = "164.26"
$ws.Cells.Item(26, 5).Value = "  +0.37%  "
$ws.Cells.Item(27, 2).Value = "Stellar"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(27, 4).Value = "0.1487"
$ws.Cells.Item(27, 5).Value = "  +54.41%  "
$ws.Cells.Item(28, 4).Value = "19.97"
$ws.Cells.Item(28, 5).Value = "  +6.87%  "
$ws.Cells.Item(29, 4).Value = "2.189"
$ws.Cells.Item(29, 5).Value = "  +16.70%  "
$ws.Cells.Item(30, 4).Value = "1.569"
$ws.Cells.Item(30, 5).Value = "  +5.64%  "
$ws.Cells.Item(31, 4).Value = "1.352"
$ws.Cells.Item(31, 5).Value = "  +2.32%  "
$ws.Cells.Item(32, 4).Value = "4.573"
$ws.Cells.Item(32, 5).Value = "  +8.49%  "
$ws.Cells.Item(33, 4).Value = "4.315"
$ws.Cells.Item(33, 5).Value = "  +4.97%  "
$ws.Cells.Item(34, 4).Value = "0.05163"
$ws.Cells.Item(34, 5).Value = "  +7.46%  "
$ws.Cells.Item(35, 4).Value = "1.216"
$ws.Cells.Item(35, 5).Value = "  +8.61%  "
$ws.Cells.Item(36, 4).Value = "0.7582"
$ws.Cells.Item(36, 5).Value = "  +10.27%  "
$ws.Cells.Item(37, 4).Value = "2.768"
$ws.Cells.Item(37, 5).Value = "  +2.14%  "
$ws.Cells.Item(38, 4).Value = "0.02002"
$ws.Cells.Item(38, 5).Value = "  +6.39%  "
$ws.Cells.Item(39, 4).Value = "2.904"
$ws.Cells.Item(39, 5).Value = "  +3.35%  "
$ws.Cells.Item(40, 4).Value = "6.591"
$ws.Cells.Item(40, 5).Value = "  +5.94%  "
$ws.Cells.Item(41, 4).Value = "78.14"
$ws.Cells.Item(41, 5).Value = "  +5.36%  "
$ws.Cells.Item(42, 4).Value = "0.4685"
$ws.Cells.Item(42, 5).Value = "  +10.99%  "
$ws.Cells.Item(43, 4).Value = "2.071"
$ws.Cells.Item(43, 5).Value = "  +7.21%  "
$ws.Cells.Item(44, 4).Value = "0.8510"
$ws.Cells.Item(44, 5).Value = "  +3.55%  "
$ws.Cells.Item(45, 4).Value = "104.59"
$ws.Cells.Item(45, 5).Value = "  +3.69%  "
$ws.Cells.Item(46, 4).Value = "0.9979"
$ws.Cells.Item(46, 5).Value = "  -0.15%  "
$ws.Cells.Item(47, 4).Value = "9.974"
$ws.Cells.Item(47, 5).Value = "  +3.87%  "
$ws.Cells.Item(48, 4).Value = "7.500"
$ws.Cells.Item(48, 5).Value = "  +7.82%  "
$ws.Cells.Item(49, 2).Value = "Decentraland"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(49, 4).Value = "0.4302"
$ws.Cells.Item(49, 5).Value = "  +9.89%  "
$ws.Cells.Item(50, 2).Value = "Elrond"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(50, 4).Value = "36.65"
$ws.Cells.Item(50, 5).Value = "  +4.00%  "
$ws.Cells.Item(51, 4).Value = "0.1196"
$ws.Cells.Item(51, 5).Value = "  +12.87%  "

$dataRange.ClearFormats()
